# Appends new listings and refreshes the "取得日時" timestamp for the
# "ランサーズ" sheet, matching the 2025-11-25 18:28 JST scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks/content in the data region so the sheet can be
# rebuilt cleanly (header row 1 is left untouched).
$ws.Hyperlinks.Delete()
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -ge 2) {
    $ws.Range("A2:H" + $lastRow).Clear()
}

$rows = @(
    @{A='2025-11-25 18:28:10'; B='【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'; C='システム開発'; D='50,000 円 ~ 100,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5405023'; G=178; H='★bot ◆ツール'}
    @{A='2025-11-25 18:28:10'; B='【急募】縫製工場向けPL・CF可視化アプリのMVP開発'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5440957'; G=93; H='◆開発 ◇アプリ'}
    @{A='2025-11-25 18:28:10'; B='【長期募集】クラウドサービス開発・保守エンジニアを求む!'; C='システム開発'; D='1,000 ~ 5,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5440461'; G=60; H='◆開発'}
    @{A='2025-11-25 18:28:10'; B='【急募】古いPHPとPerlプログラムのアップデート依頼'; C='システム開発'; D='100,000 円 ~ 200,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5440861'; G=33; H='○PHP'}
    @{A='2025-11-25 18:28:10'; B='【急募】シティヘブンの出勤情報を自動取得・管理したい!'; C='システム開発'; D='20,000 円 ~ 50,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5440436'; G=33; H='◇管理'}
    @{A='2025-11-25 18:28:10'; B='簡易サイト修正・その他小規模タスク依頼'; C='システム開発'; D='10,000 円 ~ 20,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5441146'; G=30; H='◇サイト'}
    @{A='2025-11-25 18:28:10'; B='進行管理およびチームディレクションを担当'; C='システム開発'; D='~ 5,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5418064'; G=30; H='◇管理'}
    @{A='2025-11-25 18:28:10'; B='初回 n8n+Gemini+Typefully+GoogleスプレッドのX/Threads自動投稿システム'; C='システム開発'; D='100,000 円 ~ 200,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5440440'; G=33; H=$null}
    @{A='2025-11-25 18:28:10'; B='急募 限定公開 PR 限定公開の仕事'; C='システム開発'; D='300,000 円 ~ 500,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5440230'; G=25; H=$null}
    @{A='2025-11-25 18:28:10'; B='〖リモート可〗Delphiエンジニア募集'; C='システム開発'; D='300,000 円 ~ 500,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5341051'; G=25; H=$null}
    @{A='2025-11-25 18:28:10'; B='【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5435080'; G=25; H=$null}
    @{A='2025-11-25 18:28:10'; B='【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5435079'; G=25; H=$null}
    @{A='2025-11-25 18:28:10'; B='当組織のエンジニア追加募集。'; C='システム開発'; D='100,000 円 ~ 200,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5441084'; G=18; H=$null}
    @{A='2025-11-25 18:28:10'; B='【急募】弊社Websiteの保守運用をお任せできる方を探しています!'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5440806'; G=18; H=$null}
    @{A='2025-11-25 18:28:10'; B='【急募】n8n ワークフロー実装とGoogle Sheets作成(70万円 ~ )'; C='システム開発'; D='1,000 ~ 5,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5441082'; G=10; H=$null}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row.F, $null, $null, $row.F) | Out-Null
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
    $r = $r + 1
}

